$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
@(2, "atcoder_AGC007_C.py", "Compilation Error"),
@(3, "atcoder_ABC178_A.py", "Compilation Error"),
@(4, "codeforces_203_A.py", "Compilation Error"),
@(5, "codeforces_579_A.py", "Compilation Error"),
@(6, "atcoder_AGC025_A.py", "Compilation Error"),
@(7, "atcoder_ABC169_C.py", "Runtime Error"),
@(8, "atcoder_ABC174_C.py", "Runtime Error"),
@(9, "atcoder_ABC051_A.py", "Runtime Error"),
@(10, "atcoder_ABC132_F.py", "Runtime Error"),
@(11, "codeforces_678_A.py", "Runtime Error"),
@(12, "codeforces_334_A.py", "Runtime Error"),
@(13, "codeforces_189_A.py", "Runtime Error"),
@(14, "atcoder_ABC124_C.py", "Runtime Error"),
@(15, "codeforces_96_B.py", "Runtime Error"),
@(16, "codeforces_306_A.py", "Runtime Error"),
@(17, "atcoder_ABC122_D.py", "Test Failed"),
@(18, "atcoder_ABC108_B.py", "Test Failed"),
@(19, "codeforces_651_A.py", "Test Failed"),
@(20, "codeforces_171_A.py", "Test Failed"),
@(21, "atcoder_AGC046_A.py", "Test Failed"),
@(22, "atcoder_ARC102_C.py", "Test Failed"),
@(23, "codeforces_276_B.py", "Test Failed"),
@(24, "codeforces_86_A.py", "Test Failed"),
@(25, "codeforces_544_B.py", "Test Failed"),
@(26, "atcoder_ABC042_A.py", "Test Failed"),
@(27, "atcoder_ABC136_B.py", "Test Failed"),
@(28, "codeforces_581_A.py", "Test Failed"),
@(29, "codeforces_242_A.py", "Test Failed"),
@(30, "codeforces_569_A.py", "Test Failed"),
@(31, "atcoder_ABC070_B.py", "Test Failed"),
@(32, "codeforces_672_A.py", "Test Failed"),
@(33, "codeforces_55_A.py", "Test Failed"),
@(34, "codeforces_373_B.py", "Test Failed"),
@(35, "atcoder_ABC043_B.py", "Test Failed"),
@(36, "codeforces_678_B.py", "Test Failed"),
@(37, "atcoder_ABC120_C.py", "Test Failed"),
@(38, "codeforces_379_A.py", "Test Failed"),
@(39, "codeforces_79_A.py", "Test Failed"),
@(40, "atcoder_ABC143_A.py", "Test Failed"),
@(41, "atcoder_ABC153_A.py", "Test Failed"),
@(42, "codeforces_110_B.py", "Test Failed"),
@(43, "codeforces_459_A.py", "Test Failed"),
@(44, "atcoder_ABC178_B.py", "Test Failed"),
@(45, "atcoder_AGC046_B.py", "Test Failed"),
@(46, "atcoder_ABC172_D.py", "Test Failed"),
@(47, "codeforces_369_B.py", "Test Failed"),
@(48, "codeforces_546_A.py", "Test Failed"),
@(49, "atcoder_ABC149_B.py", "Test Failed"),
@(50, "atcoder_ABC151_A.py", "Test Failed"),
@(51, "atcoder_ABC168_C.py", "Test Failed"),
@(52, "codeforces_58_B.py", "Test Failed"),
@(53, "atcoder_ABC132_A.py", "Test Failed"),
@(54, "codeforces_99_A.py", "Test Failed"),
@(55, "codeforces_514_A.py", "Test Failed"),
@(56, "codeforces_8_B.py", "Test Failed"),
@(57, "atcoder_ABC127_B.py", "Test Failed"),
@(58, "atcoder_ABC158_B.py", "Test Failed"),
@(59, "codeforces_92_A.py", "Test Failed"),
@(60, "atcoder_ABC124_A.py", "Test Failed"),
@(61, "codeforces_59_A.py", "Test Failed"),
@(62, "codeforces_49_A.py", "Test Failed"),
@(63, "atcoder_ABC149_C.py", "Test Failed"),
@(64, "atcoder_ABC139_B.py", "Test Failed"),
@(65, "codeforces_30_A.py", "Test Failed"),
@(66, "codeforces_622_A.py", "Test Failed"),
@(67, "atcoder_ABC164_A.py", "Test Failed"),
@(68, "atcoder_AGC002_A.py", "Test Failed"),
@(69, "atcoder_ABC158_A.py", "Test Failed"),
@(70, "atcoder_ARC062_B.py", "Test Failed"),
@(71, "atcoder_ABC155_E.py", "Test Failed"),
@(72, "atcoder_ABC142_A.py", "Test Failed"),
@(73, "atcoder_ABC170_A.py", "Test Failed"),
@(74, "codeforces_190_A.py", "Test Failed"),
@(75, "atcoder_ABC114_C.py", "Test Failed"),
@(76, "atcoder_ABC125_A.py", "Test Failed"),
@(77, "atcoder_ABC169_D.py", "Test Failed"),
@(78, "codeforces_32_B.py", "Test Failed"),
@(79, "codeforces_669_A.py", "Infinite Loop"),
@(80, "codeforces_340_A.py", "Infinite Loop"),
@(81, "codeforces_147_A.py", "Correct"),
@(82, "atcoder_AGC006_B.py", "Correct")
)

foreach ($row in $rows) {
    $r = $row[0]
    $filename = $row[1]
    $impact = $row[2]
    $ws.Cells.Item($r, 4).Value = $filename
    $ws.Cells.Item($r, 5).Value = $impact
}
